{"js": "// Applies the \"Cau 5\" addition + Normal-style line-spacing tweak described\n// by the diff. Body is `async (context) => { ... }`.\n\n// ---------------------------------------------------------------------\n// 1) Find the last non-empty paragraph in the body (\"S\u1ebd b\u1ecb d\u00ednh xung \u0111\u1ed9t\n//    ngay kh\u00fac n\u00e0y.\") and append the new \"C\u00e2u 5\" block right after it,\n//    ahead of the trailing blank paragraphs that precede the sectPr.\n// ---------------------------------------------------------------------\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nlet anchorIndex = -1;\nfor (let i = items.length - 1; i >= 0; i--) {\n  if (items[i].text.trim().length > 0) {\n    anchorIndex = i;\n    break;\n  }\n}\nif (anchorIndex === -1) {\n  anchorIndex = items.length - 1;\n}\n\nconst newParagraphTexts = [\n  \"C\u00e2u 5.\",\n  \"- H\u00e3y m\u1edf file world trong c\u00e2u 1 trong nh\u00e1nh M\u00e3 s\u1ed1 SV c\u1ee7 b\u1ea1n \u0111\u1ec3 th\u00eam ph\u1ea7n h\u01b0\u1edbng d\u1eabn ch\u1ea1y 2 file tr\u00ean.\",\n  \"\u0110\u1ed1i v\u1edbi file.cpp\",\n  \"M\u1edf visual code c++ new open file v\u00e0 ch\u1ea1y\",\n  \"\u0110\u1ed1i v\u1edbi file.sh\",\n  \"C\u1ea5p quy\u1ec1n cho file th\u1ef1c thi:\",\n  \"Chmod +x file.sh\",\n  \"Ch\u1ea1y\",\n  \"./file.sh\",\n  \"\"\n];\n\nlet cursor = items[anchorIndex];\nfor (const text of newParagraphTexts) {\n  cursor = cursor.insertParagraph(text, \"After\");\n}\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 2) Normal style: line spacing becomes 256/20 = 12.8pt with the \"auto\"\n//    (multiple) rule, i.e. <w:spacing w:line=\"256\" w:lineRule=\"auto\"/>.\n// ---------------------------------------------------------------------\nconst normalStyle = context.document.styles.getByName(\"Normal\");\nconst normalFormat = normalStyle.paragraphFormat;\nnormalFormat.lineSpacing = 12.8;\n// `lineSpacingRule` isn't surfaced as a public Word.ParagraphFormat member\n// of this shim, but the OM bridge understands it (same bridge Word COM's\n// `ParagraphFormat.LineSpacingRule = wdLineSpaceMultiple` drives) \u2014 reach\n// it directly so the saved style carries `w:lineRule=\"auto\"` too.\nnormalFormat._omSet(\"LineSpacingRule\", 5 /* wdLineSpaceMultiple */);\nawait context.sync();\n", "ps1": "# Applies the \"Cau 5\" addition + Normal-style line-spacing tweak described\n# by the diff. $d / $word / $app are pre-seeded by the host; $d is the\n# ActiveDocument.\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# 1) Find the last non-empty paragraph in the body (\"S\u1ebd b\u1ecb d\u00ednh xung \u0111\u1ed9t\n#    ngay kh\u00fac n\u00e0y.\") and append the new \"C\u00e2u 5\" block right after it,\n#    ahead of the trailing blank paragraphs that precede the sectPr.\n# ---------------------------------------------------------------------\n$count = $d.Paragraphs.Count\n$anchorIndex = -1\nfor ($i = $count; $i -ge 1; $i--) {\n  $t = $d.Paragraphs($i).Range.Text\n  if ($t.Trim().Length -gt 0) {\n    $anchorIndex = $i\n    break\n  }\n}\nif ($anchorIndex -eq -1) {\n  $anchorIndex = $count\n}\n\n$newParagraphTexts = @(\n  \"C\u00e2u 5.\",\n  \"- H\u00e3y m\u1edf file world trong c\u00e2u 1 trong nh\u00e1nh M\u00e3 s\u1ed1 SV c\u1ee7 b\u1ea1n \u0111\u1ec3 th\u00eam ph\u1ea7n h\u01b0\u1edbng d\u1eabn ch\u1ea1y 2 file tr\u00ean.\",\n  \"\u0110\u1ed1i v\u1edbi file.cpp\",\n  \"M\u1edf visual code c++ new open file v\u00e0 ch\u1ea1y\",\n  \"\u0110\u1ed1i v\u1edbi file.sh\",\n  \"C\u1ea5p quy\u1ec1n cho file th\u1ef1c thi:\",\n  \"Chmod +x file.sh\",\n  \"Ch\u1ea1y\",\n  \"./file.sh\",\n  \"\"\n)\n\n$idx = $anchorIndex\nforeach ($text in $newParagraphTexts) {\n  $d.Paragraphs($idx).Range.InsertParagraphAfter()\n  $idx = $idx + 1\n  if ($text.Length -gt 0) {\n    $d.Paragraphs($idx).Range.InsertAfter($text)\n  }\n}\n\n# ---------------------------------------------------------------------\n# 2) Normal style: line spacing becomes 256/20 = 12.8pt with the \"auto\"\n#    (multiple) rule, i.e. <w:spacing w:line=\"256\" w:lineRule=\"auto\"/>.\n# ---------------------------------------------------------------------\n$normalStyle = $d.Styles(\"Normal\")\n$normalStyle.ParagraphFormat.LineSpacingRule = 5  # wdLineSpaceMultiple\n$normalStyle.ParagraphFormat.LineSpacing = 12.8\n"}
